$d = $word.ActiveDocument

# Step 1: remove first paragraph "NLP Project: Paperwork"
$d.Paragraphs.Item(1).Range.Delete()

# Step 2: update heading formatting sz 28->32 (14->16 pts)
$heading = $d.Paragraphs.Item(1)
$heading.Range.Font.Size = 16

# Step 3: replace paragraph 2 (first body para) text
$p2 = $d.Paragraphs.Item(2)
$p2r = $d.Range($p2.Range.Start, $p2.Range.End - 1)
$p2r.Text = "This section presents the outcome in their preliminary experiments."

# Step 4: replace paragraph 3 (second body para) text
$p3 = $d.Paragraphs.Item(3)
$p3r = $d.Range($p3.Range.Start, $p3.Range.End - 1)
$p3r.Text = "The data was acquired from social networking sites such as Facebook, Twitter and YouTube. Moreover, two classification tasks were carried out: cyberbullying event detection and the classification of text categories related to cyberbullying. Using the Special Text Replacement function in Excel, the dataset was normalized. The researchers get the frequency of each instance of the word under each category. Out of 625 statements that were extracted from the social media sites, the harmfulness level 1 was most prevalent with a frequency of 33%. It was followed by non-cyberbullying events, with a frequency of 41%. The occurrences of severe cyberbullying events were least prevalent among the three classification; it has a frequency of 26%."

# Step 5
$p4 = $d.Paragraphs.Item(4)
$pos = $p4.Range.Start
Write-Output "pos=$pos"
$bmX = $d.Bookmarks.Item("_GoBack")
Write-Output "bookmark before insert: start=$($bmX.Range.Start) end=$($bmX.Range.End)"

$beforeText = "At the second level of annotation, the bad description was the most prevalent with a frequency of 27%, it was followed by Social Rejection (18%), Intelligence (16%), Sexuality (13%), Physical Appearance (10%) and the category Race an"
$insertRange = $d.Range($pos, $pos)
$insertRange.InsertBefore($beforeText)

$bmY = $d.Bookmarks.Item("_GoBack")
Write-Output "bookmark after insertBefore: start=$($bmY.Range.Start) end=$($bmY.Range.End)"
Write-Output "expected = $($pos + $beforeText.Length)"

$fmt1 = $d.Range($pos, $pos + $beforeText.Length)
$fmt1.Font.Name = "Arial"
$fmt1.Font.Size = 12

$bmZ = $d.Bookmarks.Item("_GoBack")
Write-Output "bookmark after fmt1: start=$($bmZ.Range.Start) end=$($bmZ.Range.End)"

$afterText = "d Culture was the least prevalent among the six categories, it has a frequency of 8%."
$bm = $d.Bookmarks.Item("_GoBack")
$afterPos = $bm.Range.End
Write-Output "afterPos=$afterPos"
$insertRange2 = $d.Range($afterPos, $afterPos)
$insertRange2.InsertAfter($afterText)

$p4text = $d.Paragraphs.Item(4).Range.Text
Write-Output "p4 text: [$p4text]"

$bmFinal = $d.Bookmarks.Item("_GoBack")
Write-Output "bookmark final: start=$($bmFinal.Range.Start) end=$($bmFinal.Range.End)"

Write-Output "char at 1067-1069: [$($d.Range(1060,1075).Text)]"
Write-Output "char at bookmark pos +/- 5: [$($d.Range($bmFinal.Range.Start - 5, $bmFinal.Range.Start + 5).Text)]"
